# Commit message: "Changed Q3 to Q4 in spreadsheets"
# Replace every occurrence of "Q3" with "Q4" across all cell values in the
# workbook (mirrors an Excel Find & Replace > Replace All run against the
# whole workbook).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    if ($used -ne $null) {
        [void]$used.Replace("Q3", "Q4")
    }
}

# The "file" sheet was the active tab while editing and ended up with the
# cursor on I18 afterwards.
$fileSheet = $wb.Worksheets.Item("file")
[void]$fileSheet.Activate()
[void]$fileSheet.Range("I18").Select()
